# Vai ao fim do documento e acrescenta as novas linhas, tal como um
# utilizador faria ao escrever no Word: Enter (linha em branco), Enter,
# "Adicionei conteúdo ok", Enter, "E volto a adicionar conteúdo".
$d   = $word.ActiveDocument
$sel = $word.Selection

[void]$sel.EndKey(6)            # wdStory -> vai para o fim do documento
[void]$sel.TypeParagraph()      # fecha o parágrafo "Vou contar uma historia"

# Marca temporariamente a nova linha em branco para o motor criar o novo
# parágrafo e depois limpa-a, ficando um <w:p/> completamente vazio (tal
# como aconteceria com uma linha em branco real, sem qualquer carácter).
[void]$sel.TypeText("~")
[void]$sel.TypeParagraph()
[void]$sel.TypeText("Adicionei conteúdo ok")
[void]$sel.TypeParagraph()
[void]$sel.TypeText("E volto a adicionar conteúdo")

$blankPara = $d.Paragraphs(2).Range
$blankPara.MoveEnd(1, -1)       # wdCharacter, exclui a marca de parágrafo
$blankPara.Text = ""
